$wb = $excel.ActiveWorkbook

# 1) Add the new "PO Forecast" sheet at the end, matching sheetId=3 / rId3
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$poForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$poForecast.Name = "PO Forecast"

$weeklySheet = $wb.Worksheets.Item(1)
$monthlySheet = $wb.Worksheets.Item(2)

# 2) Rename the header cells on the existing sheets
$weeklySheet.Range("B1").Value = "Weekly_PO_Qty"
$monthlySheet.Range("B1").Value = "Monthly_PO_Qty"

# 3) Copy header formatting (bold/border/centered style) from an existing header row
$weeklySheet.Range("A1:B1").Copy()
$poForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-formatted style (numFmt 165) from an existing date cell
$weeklySheet.Range("A2").Copy()
$poForecast.Range("A2:A77").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Write header values
$poForecast.Range("A1").Value = "ds"
$poForecast.Range("B1").Value = "PO_Forecast"
$poForecast.Range("C1").Value = "yhat_lower"
$poForecast.Range("D1").Value = "yhat_upper"

# 5) Write the forecast data rows
$poForecast.Cells.Item(2, 1).Value = 44934.99999999999
$poForecast.Cells.Item(2, 2).Value = 536
$poForecast.Cells.Item(2, 3).Value = -451.1103726907901
$poForecast.Cells.Item(2, 4).Value = 1484.190099723694
$poForecast.Cells.Item(3, 1).Value = 44941.99999999999
$poForecast.Cells.Item(3, 2).Value = 541
$poForecast.Cells.Item(3, 3).Value = -381.2405740727347
$poForecast.Cells.Item(3, 4).Value = 1489.360270606441
$poForecast.Cells.Item(4, 1).Value = 44948.99999999999
$poForecast.Cells.Item(4, 2).Value = 547
$poForecast.Cells.Item(4, 3).Value = -393.8449651314277
$poForecast.Cells.Item(4, 4).Value = 1480.160776447115
$poForecast.Cells.Item(5, 1).Value = 44955.99999999999
$poForecast.Cells.Item(5, 2).Value = 552
$poForecast.Cells.Item(5, 3).Value = -379.1415381900345
$poForecast.Cells.Item(5, 4).Value = 1525.609952188099
$poForecast.Cells.Item(6, 1).Value = 44962.99999999999
$poForecast.Cells.Item(6, 2).Value = 557
$poForecast.Cells.Item(6, 3).Value = -436.2775986683242
$poForecast.Cells.Item(6, 4).Value = 1465.50088364395
$poForecast.Cells.Item(7, 1).Value = 44990.99999999999
$poForecast.Cells.Item(7, 2).Value = 578
$poForecast.Cells.Item(7, 3).Value = -370.3286166855503
$poForecast.Cells.Item(7, 4).Value = 1493.272965508834
$poForecast.Cells.Item(8, 1).Value = 44997.99999999999
$poForecast.Cells.Item(8, 2).Value = 584
$poForecast.Cells.Item(8, 3).Value = -341.557430860952
$poForecast.Cells.Item(8, 4).Value = 1571.089553341598
$poForecast.Cells.Item(9, 1).Value = 45004.99999999999
$poForecast.Cells.Item(9, 2).Value = 589
$poForecast.Cells.Item(9, 3).Value = -388.2891186994665
$poForecast.Cells.Item(9, 4).Value = 1552.210421155645
$poForecast.Cells.Item(10, 1).Value = 45011.99999999999
$poForecast.Cells.Item(10, 2).Value = 594
$poForecast.Cells.Item(10, 3).Value = -338.5450057348278
$poForecast.Cells.Item(10, 4).Value = 1639.422945497007
$poForecast.Cells.Item(11, 1).Value = 45025.99999999999
$poForecast.Cells.Item(11, 2).Value = 605
$poForecast.Cells.Item(11, 3).Value = -278.4334386220742
$poForecast.Cells.Item(11, 4).Value = 1554.738365922371
$poForecast.Cells.Item(12, 1).Value = 45039.99999999999
$poForecast.Cells.Item(12, 2).Value = 615
$poForecast.Cells.Item(12, 3).Value = -244.8635608251109
$poForecast.Cells.Item(12, 4).Value = 1586.417382905146
$poForecast.Cells.Item(13, 1).Value = 45046.99999999999
$poForecast.Cells.Item(13, 2).Value = 621
$poForecast.Cells.Item(13, 3).Value = -383.3981447075282
$poForecast.Cells.Item(13, 4).Value = 1553.575801290722
$poForecast.Cells.Item(14, 1).Value = 45067.99999999999
$poForecast.Cells.Item(14, 2).Value = 637
$poForecast.Cells.Item(14, 3).Value = -274.1893057694791
$poForecast.Cells.Item(14, 4).Value = 1573.11083181097
$poForecast.Cells.Item(15, 1).Value = 45081.99999999999
$poForecast.Cells.Item(15, 2).Value = 647
$poForecast.Cells.Item(15, 3).Value = -285.9872741922529
$poForecast.Cells.Item(15, 4).Value = 1582.201448089284
$poForecast.Cells.Item(16, 1).Value = 45088.99999999999
$poForecast.Cells.Item(16, 2).Value = 652
$poForecast.Cells.Item(16, 3).Value = -312.7924631985493
$poForecast.Cells.Item(16, 4).Value = 1599.346745236295
$poForecast.Cells.Item(17, 1).Value = 45095.99999999999
$poForecast.Cells.Item(17, 2).Value = 658
$poForecast.Cells.Item(17, 3).Value = -302.4125431402269
$poForecast.Cells.Item(17, 4).Value = 1591.8699398821
$poForecast.Cells.Item(18, 1).Value = 45102.99999999999
$poForecast.Cells.Item(18, 2).Value = 663
$poForecast.Cells.Item(18, 3).Value = -297.6333661396986
$poForecast.Cells.Item(18, 4).Value = 1576.121311201682
$poForecast.Cells.Item(19, 1).Value = 45109.99999999999
$poForecast.Cells.Item(19, 2).Value = 668
$poForecast.Cells.Item(19, 3).Value = -408.8206693047049
$poForecast.Cells.Item(19, 4).Value = 1615.246307599122
$poForecast.Cells.Item(20, 1).Value = 45116.99999999999
$poForecast.Cells.Item(20, 2).Value = 674
$poForecast.Cells.Item(20, 3).Value = -380.5973006427542
$poForecast.Cells.Item(20, 4).Value = 1601.906781226941
$poForecast.Cells.Item(21, 1).Value = 45123.99999999999
$poForecast.Cells.Item(21, 2).Value = 679
$poForecast.Cells.Item(21, 3).Value = -309.1177047435612
$poForecast.Cells.Item(21, 4).Value = 1693.2958702026
$poForecast.Cells.Item(22, 1).Value = 45130.99999999999
$poForecast.Cells.Item(22, 2).Value = 684
$poForecast.Cells.Item(22, 3).Value = -261.5576796538203
$poForecast.Cells.Item(22, 4).Value = 1576.205061085093
$poForecast.Cells.Item(23, 1).Value = 45137.99999999999
$poForecast.Cells.Item(23, 2).Value = 689
$poForecast.Cells.Item(23, 3).Value = -257.9903904359527
$poForecast.Cells.Item(23, 4).Value = 1612.587060756734
$poForecast.Cells.Item(24, 1).Value = 45144.99999999999
$poForecast.Cells.Item(24, 2).Value = 695
$poForecast.Cells.Item(24, 3).Value = -274.4756466199734
$poForecast.Cells.Item(24, 4).Value = 1678.530899128802
$poForecast.Cells.Item(25, 1).Value = 45151.99999999999
$poForecast.Cells.Item(25, 2).Value = 700
$poForecast.Cells.Item(25, 3).Value = -238.5368743858807
$poForecast.Cells.Item(25, 4).Value = 1597.901861656718
$poForecast.Cells.Item(26, 1).Value = 45158.99999999999
$poForecast.Cells.Item(26, 2).Value = 705
$poForecast.Cells.Item(26, 3).Value = -186.6031932482897
$poForecast.Cells.Item(26, 4).Value = 1684.167300458852
$poForecast.Cells.Item(27, 1).Value = 45165.99999999999
$poForecast.Cells.Item(27, 2).Value = 711
$poForecast.Cells.Item(27, 3).Value = -244.9434629423955
$poForecast.Cells.Item(27, 4).Value = 1716.055345115474
$poForecast.Cells.Item(28, 1).Value = 45179.99999999999
$poForecast.Cells.Item(28, 2).Value = 721
$poForecast.Cells.Item(28, 3).Value = -250.1194911851618
$poForecast.Cells.Item(28, 4).Value = 1690.362382499033
$poForecast.Cells.Item(29, 1).Value = 45186.99999999999
$poForecast.Cells.Item(29, 2).Value = 727
$poForecast.Cells.Item(29, 3).Value = -261.6492281899672
$poForecast.Cells.Item(29, 4).Value = 1671.497549507202
$poForecast.Cells.Item(30, 1).Value = 45193.99999999999
$poForecast.Cells.Item(30, 2).Value = 732
$poForecast.Cells.Item(30, 3).Value = -247.8994083198549
$poForecast.Cells.Item(30, 4).Value = 1731.820713169017
$poForecast.Cells.Item(31, 1).Value = 45200.99999999999
$poForecast.Cells.Item(31, 2).Value = 737
$poForecast.Cells.Item(31, 3).Value = -230.8991548770761
$poForecast.Cells.Item(31, 4).Value = 1671.302409139538
$poForecast.Cells.Item(32, 1).Value = 45207.99999999999
$poForecast.Cells.Item(32, 2).Value = 742
$poForecast.Cells.Item(32, 3).Value = -235.2708750832266
$poForecast.Cells.Item(32, 4).Value = 1706.318975063748
$poForecast.Cells.Item(33, 1).Value = 45214.99999999999
$poForecast.Cells.Item(33, 2).Value = 748
$poForecast.Cells.Item(33, 3).Value = -205.8685595202358
$poForecast.Cells.Item(33, 4).Value = 1686.976265907868
$poForecast.Cells.Item(34, 1).Value = 45221.99999999999
$poForecast.Cells.Item(34, 2).Value = 753
$poForecast.Cells.Item(34, 3).Value = -189.8356705124136
$poForecast.Cells.Item(34, 4).Value = 1675.8021668613
$poForecast.Cells.Item(35, 1).Value = 45228.99999999999
$poForecast.Cells.Item(35, 2).Value = 758
$poForecast.Cells.Item(35, 3).Value = -217.8506235436414
$poForecast.Cells.Item(35, 4).Value = 1717.155000602488
$poForecast.Cells.Item(36, 1).Value = 45235.99999999999
$poForecast.Cells.Item(36, 2).Value = 764
$poForecast.Cells.Item(36, 3).Value = -220.4401278479138
$poForecast.Cells.Item(36, 4).Value = 1725.342529417069
$poForecast.Cells.Item(37, 1).Value = 45242.99999999999
$poForecast.Cells.Item(37, 2).Value = 769
$poForecast.Cells.Item(37, 3).Value = -187.3582541065503
$poForecast.Cells.Item(37, 4).Value = 1717.033731888069
$poForecast.Cells.Item(38, 1).Value = 45249.99999999999
$poForecast.Cells.Item(38, 2).Value = 774
$poForecast.Cells.Item(38, 3).Value = -182.6226299476634
$poForecast.Cells.Item(38, 4).Value = 1686.992568217247
$poForecast.Cells.Item(39, 1).Value = 45256.99999999999
$poForecast.Cells.Item(39, 2).Value = 779
$poForecast.Cells.Item(39, 3).Value = -110.3174756697013
$poForecast.Cells.Item(39, 4).Value = 1735.177320027163
$poForecast.Cells.Item(40, 1).Value = 45263.99999999999
$poForecast.Cells.Item(40, 2).Value = 785
$poForecast.Cells.Item(40, 3).Value = -133.0131600809226
$poForecast.Cells.Item(40, 4).Value = 1706.24792026145
$poForecast.Cells.Item(41, 1).Value = 45277.99999999999
$poForecast.Cells.Item(41, 2).Value = 795
$poForecast.Cells.Item(41, 3).Value = -189.3895146727392
$poForecast.Cells.Item(41, 4).Value = 1746.526264791298
$poForecast.Cells.Item(42, 1).Value = 45298.99999999999
$poForecast.Cells.Item(42, 2).Value = 811
$poForecast.Cells.Item(42, 3).Value = -137.9537087611872
$poForecast.Cells.Item(42, 4).Value = 1765.291510620813
$poForecast.Cells.Item(43, 1).Value = 45305.99999999999
$poForecast.Cells.Item(43, 2).Value = 817
$poForecast.Cells.Item(43, 3).Value = -74.54186150942688
$poForecast.Cells.Item(43, 4).Value = 1735.972195566829
$poForecast.Cells.Item(44, 1).Value = 45312.99999999999
$poForecast.Cells.Item(44, 2).Value = 822
$poForecast.Cells.Item(44, 3).Value = -132.4362693125459
$poForecast.Cells.Item(44, 4).Value = 1733.254191000826
$poForecast.Cells.Item(45, 1).Value = 45326.99999999999
$poForecast.Cells.Item(45, 2).Value = 832
$poForecast.Cells.Item(45, 3).Value = -169.6114836548221
$poForecast.Cells.Item(45, 4).Value = 1819.257267941574
$poForecast.Cells.Item(46, 1).Value = 45333.99999999999
$poForecast.Cells.Item(46, 2).Value = 838
$poForecast.Cells.Item(46, 3).Value = -90.50440758460685
$poForecast.Cells.Item(46, 4).Value = 1819.112168228465
$poForecast.Cells.Item(47, 1).Value = 45403.99999999999
$poForecast.Cells.Item(47, 2).Value = 891
$poForecast.Cells.Item(47, 3).Value = -17.88467730007453
$poForecast.Cells.Item(47, 4).Value = 1791.015204558234
$poForecast.Cells.Item(48, 1).Value = 45410.99999999999
$poForecast.Cells.Item(48, 2).Value = 896
$poForecast.Cells.Item(48, 3).Value = -30.60439914063287
$poForecast.Cells.Item(48, 4).Value = 1891.15187501829
$poForecast.Cells.Item(49, 1).Value = 45417.99999999999
$poForecast.Cells.Item(49, 2).Value = 901
$poForecast.Cells.Item(49, 3).Value = -54.22823494102234
$poForecast.Cells.Item(49, 4).Value = 1854.01464492524
$poForecast.Cells.Item(50, 1).Value = 45424.99999999999
$poForecast.Cells.Item(50, 2).Value = 907
$poForecast.Cells.Item(50, 3).Value = -32.40733419614841
$poForecast.Cells.Item(50, 4).Value = 1799.465560947466
$poForecast.Cells.Item(51, 1).Value = 45431.99999999999
$poForecast.Cells.Item(51, 2).Value = 912
$poForecast.Cells.Item(51, 3).Value = -66.42096398270839
$poForecast.Cells.Item(51, 4).Value = 1804.233618886565
$poForecast.Cells.Item(52, 1).Value = 45438.99999999999
$poForecast.Cells.Item(52, 2).Value = 917
$poForecast.Cells.Item(52, 3).Value = -27.64821714416603
$poForecast.Cells.Item(52, 4).Value = 1839.675851974034
$poForecast.Cells.Item(53, 1).Value = 45445.99999999999
$poForecast.Cells.Item(53, 2).Value = 922
$poForecast.Cells.Item(53, 3).Value = -11.77158366526069
$poForecast.Cells.Item(53, 4).Value = 1891.556900957179
$poForecast.Cells.Item(54, 1).Value = 45452.99999999999
$poForecast.Cells.Item(54, 2).Value = 928
$poForecast.Cells.Item(54, 3).Value = -19.15621604466707
$poForecast.Cells.Item(54, 4).Value = 1920.430266579096
$poForecast.Cells.Item(55, 1).Value = 45459.99999999999
$poForecast.Cells.Item(55, 2).Value = 933
$poForecast.Cells.Item(55, 3).Value = -35.39926555053554
$poForecast.Cells.Item(55, 4).Value = 1918.611405728718
$poForecast.Cells.Item(56, 1).Value = 45466.99999999999
$poForecast.Cells.Item(56, 2).Value = 938
$poForecast.Cells.Item(56, 3).Value = -2.833850978264937
$poForecast.Cells.Item(56, 4).Value = 1874.290714475989
$poForecast.Cells.Item(57, 1).Value = 45473.99999999999
$poForecast.Cells.Item(57, 2).Value = 944
$poForecast.Cells.Item(57, 3).Value = 68.27519036822582
$poForecast.Cells.Item(57, 4).Value = 1914.825849704518
$poForecast.Cells.Item(58, 1).Value = 45487.99999999999
$poForecast.Cells.Item(58, 2).Value = 954
$poForecast.Cells.Item(58, 3).Value = 11.01462485420477
$poForecast.Cells.Item(58, 4).Value = 1959.735393354487
$poForecast.Cells.Item(59, 1).Value = 45515.99999999999
$poForecast.Cells.Item(59, 2).Value = 975
$poForecast.Cells.Item(59, 3).Value = 67.21794822755254
$poForecast.Cells.Item(59, 4).Value = 1918.402761915914
$poForecast.Cells.Item(60, 1).Value = 45522.99999999999
$poForecast.Cells.Item(60, 2).Value = 981
$poForecast.Cells.Item(60, 3).Value = -3.841608952173154
$poForecast.Cells.Item(60, 4).Value = 1930.428378374679
$poForecast.Cells.Item(61, 1).Value = 45529.99999999999
$poForecast.Cells.Item(61, 2).Value = 986
$poForecast.Cells.Item(61, 3).Value = 6.247980681770808
$poForecast.Cells.Item(61, 4).Value = 1907.347188865172
$poForecast.Cells.Item(62, 1).Value = 45536.99999999999
$poForecast.Cells.Item(62, 2).Value = 991
$poForecast.Cells.Item(62, 3).Value = 42.27112847326548
$poForecast.Cells.Item(62, 4).Value = 1961.146882282961
$poForecast.Cells.Item(63, 1).Value = 45557.99999999999
$poForecast.Cells.Item(63, 2).Value = 1007
$poForecast.Cells.Item(63, 3).Value = -6.846006386785565
$poForecast.Cells.Item(63, 4).Value = 1899.081594579177
$poForecast.Cells.Item(64, 1).Value = 45564.99999999999
$poForecast.Cells.Item(64, 2).Value = 1012
$poForecast.Cells.Item(64, 3).Value = 88.61192707610367
$poForecast.Cells.Item(64, 4).Value = 1975.289075863381
$poForecast.Cells.Item(65, 1).Value = 45571.99999999999
$poForecast.Cells.Item(65, 2).Value = 1018
$poForecast.Cells.Item(65, 3).Value = 46.5237702958051
$poForecast.Cells.Item(65, 4).Value = 2034.386273828083
$poForecast.Cells.Item(66, 1).Value = 45578.99999999999
$poForecast.Cells.Item(66, 2).Value = 1023
$poForecast.Cells.Item(66, 3).Value = 85.54186408073039
$poForecast.Cells.Item(66, 4).Value = 1940.56151267927
$poForecast.Cells.Item(67, 1).Value = 45585.99999999999
$poForecast.Cells.Item(67, 2).Value = 1028
$poForecast.Cells.Item(67, 3).Value = 77.0398321262694
$poForecast.Cells.Item(67, 4).Value = 1986.514649358315
$poForecast.Cells.Item(68, 1).Value = 45592.99999999999
$poForecast.Cells.Item(68, 2).Value = 1034
$poForecast.Cells.Item(68, 3).Value = 97.00220913396772
$poForecast.Cells.Item(68, 4).Value = 1993.47306249839
$poForecast.Cells.Item(69, 1).Value = 45599.99999999999
$poForecast.Cells.Item(69, 2).Value = 1039
$poForecast.Cells.Item(69, 3).Value = 31.05479167084502
$poForecast.Cells.Item(69, 4).Value = 2014.026716481661
$poForecast.Cells.Item(70, 1).Value = 45606.99999999999
$poForecast.Cells.Item(70, 2).Value = 1044
$poForecast.Cells.Item(70, 3).Value = 76.54844048084956
$poForecast.Cells.Item(70, 4).Value = 1988.538409707064
$poForecast.Cells.Item(71, 1).Value = 45613.99999999999
$poForecast.Cells.Item(71, 2).Value = 1049
$poForecast.Cells.Item(71, 3).Value = 91.57927621587417
$poForecast.Cells.Item(71, 4).Value = 1996.011157222936
$poForecast.Cells.Item(72, 1).Value = 45620.99999999999
$poForecast.Cells.Item(72, 2).Value = 1055
$poForecast.Cells.Item(72, 3).Value = 119.3961052447058
$poForecast.Cells.Item(72, 4).Value = 2028.576463457578
$poForecast.Cells.Item(73, 1).Value = 45627.99999999999
$poForecast.Cells.Item(73, 2).Value = 1060
$poForecast.Cells.Item(73, 3).Value = 88.68434852104447
$poForecast.Cells.Item(73, 4).Value = 1956.689632062236
$poForecast.Cells.Item(74, 1).Value = 45634.99999999999
$poForecast.Cells.Item(74, 2).Value = 1065
$poForecast.Cells.Item(74, 3).Value = 119.4016264357154
$poForecast.Cells.Item(74, 4).Value = 2039.425093158243
$poForecast.Cells.Item(75, 1).Value = 45641.99999999999
$poForecast.Cells.Item(75, 2).Value = 1071
$poForecast.Cells.Item(75, 3).Value = 172.7388426866139
$poForecast.Cells.Item(75, 4).Value = 2090.584417208542
$poForecast.Cells.Item(76, 1).Value = 45648.99999999999
$poForecast.Cells.Item(76, 2).Value = 1076
$poForecast.Cells.Item(76, 3).Value = 84.45752485298722
$poForecast.Cells.Item(76, 4).Value = 2092.542157659684
$poForecast.Cells.Item(77, 1).Value = 45655.99999999999
$poForecast.Cells.Item(77, 2).Value = 1081
$poForecast.Cells.Item(77, 3).Value = 76.60147940295141
$poForecast.Cells.Item(77, 4).Value = 2073.338707564757
